$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258, shifting existing rows 258..309 down to 259..310
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with its data
$ws.Cells.Item(258, 1).Value = 7
$ws.Cells.Item(258, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(258, 3).Value = "Ñuble"
$ws.Cells.Item(258, 4).Value = 44641
$ws.Cells.Item(258, 5).Value = 16
$ws.Cells.Item(258, 6).Value = 100114001
$ws.Cells.Item(258, 7).Value = "Papa"
$ws.Cells.Item(258, 8).Value = "Patagonia"
$ws.Cells.Item(258, 9).Value = "1a (cosecha)"
$ws.Cells.Item(258, 10).Value = 160
$ws.Cells.Item(258, 11).Value = 7000
$ws.Cells.Item(258, 12).Value = 7500
$ws.Cells.Item(258, 13).Value = 7250
$ws.Cells.Item(258, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(258, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(258, 16).Value = 290
$ws.Cells.Item(258, 17).Value = 25
$ws.Cells.Item(258, 18).Value = "Hortaliza"
